# Add a team win/loss/tie record to the right of the existing columns (A:AB).
# New columns: AC = Wins, AD = Losses, AE = Ties.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlTop = -4160

# Header row (row 1) — match the look of the existing header cells
# (bold font, thin box border, centered horizontally, top-aligned vertically).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$headerRange = $ws.Range("AC1:AE1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = $xlCenter
$headerRange.VerticalAlignment = $xlTop

# Every player row (2-41) gets the same team record: 77 wins, 85 losses, 0 ties.
for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 29).Value = 77   # AC
    $ws.Cells.Item($row, 30).Value = 85   # AD
    $ws.Cells.Item($row, 31).Value = 0    # AE
}
